# Sistema de Cadastro Python - update Usuarios worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# XML width = ColumnWidth (as read back via COM) + 0.83
$ws.Columns.Item(2).ColumnWidth = 15 - 0.83   # B: 10 -> 15
$ws.Columns.Item(3).ColumnWidth = 18 - 0.83   # C: 11 -> 18
$ws.Columns.Item(4).ColumnWidth = 19 - 0.83   # D: 23 -> 19
$ws.Columns.Item(6).ColumnWidth = 11 - 0.83   # F: 6  -> 11

# --- Fix typo in row 2 email ---
$ws.Range("D2").Value = "thiago@gmail.comm"

# --- Replace user in row 3 ---
$ws.Range("B3").Value = "Julia"
$ws.Range("C3").Value = "Souza"
$ws.Range("D3").Value = "Julia@gmail.com"

# --- Add new rows 4-6 ---
# Match the existing data-row look & feel (left/center aligned, same font as row 2/3)
$dataRange = $ws.Range("A4:H6")
$dataRange.HorizontalAlignment = -4131   # xlLeft
$dataRange.VerticalAlignment = -4108    # xlCenter

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "João Pedro"
$ws.Range("C4").Value = "Costa da Silva"
$ws.Range("D4").Value = "joao@gmail.com"
$ws.Range("E4").Value = 53
$ws.Range("F4").Value = "M"
$ws.Range("G4").Value = "'01234567891"
$ws.Range("H4").Value = "Aa123456789*"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Maria Eduarda"
$ws.Range("C5").Value = "da Silva Peixoto"
$ws.Range("D5").Value = "maria@gmail.com"
$ws.Range("E5").Value = 82
$ws.Range("F5").Value = "Feminino"
$ws.Range("G5").Value = "'01234567890"
$ws.Range("H5").Value = "Aa123456789*"

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Lucas"
$ws.Range("C6").Value = "alme"
$ws.Range("D6").Value = "sdg@."
$ws.Range("E6").Value = 120
$ws.Range("F6").Value = "Masculino"
$ws.Range("G6").Value = "'12301231548"
$ws.Range("H6").Value = "Aasfdsdgdj1*"
